# Apply cell-level edits described by the diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: challenges ---
$wsChallenges = $wb.Worksheets.Item("challenges")
$wsChallenges.Range("L3").Value = 8
$wsChallenges.Range("M3").Value = 3
$wsChallenges.Range("B4").Value = 3
$wsChallenges.Range("D4").Value = "G3"
$wsChallenges.Range("M4").Value = 4
$wsChallenges.Range("B5").Value = 4
$wsChallenges.Range("D5").Value = "G4"
$wsChallenges.Range("L5").Value = 1
$wsChallenges.Range("M5").Value = 5
$wsChallenges.Range("O5").Value = 3
$wsChallenges.Range("B6").Value = 5
$wsChallenges.Range("D6").Value = "G5"
$wsChallenges.Range("L6").Value = 1
$wsChallenges.Range("M6").Value = 6
$wsChallenges.Range("O6").Value = 4
$wsChallenges.Range("B7").Value = 6
$wsChallenges.Range("D7").Value = "G6"
$wsChallenges.Range("L7").Value = 6
$wsChallenges.Range("M7").ClearContents()
$wsChallenges.Range("O7").Value = 5

# Row 8 (campaign 17, id 8) was removed entirely; delete the row so the sheet dimension shrinks to A1:O7.
$wsChallenges.Rows.Item(8).Delete()

# --- Sheet: tasks ---
$wsTasks = $wb.Worksheets.Item("tasks")
$wsTasks.Range("L2").Value = " [SECRET, EQUAL, h5bdjcop3njkonq5a]"
$wsTasks.Range("B3").Value = "Take_200_steps"
$wsTasks.Range("L3").Value = "[STEPS, STRICTLY_GREATER, 200], [SECRET, EQUAL, d1vnnashgx3onjzzjzjm]"
$wsTasks.Range("B4").Value = "Take_300_steps"
$wsTasks.Range("I4").Value = "WALK"
$wsTasks.Range("J4").Value = "WALK"
$wsTasks.Range("L4").Value = "[STEPS, STRICTLY_GREATER, 300], [SECRET, EQUAL, pku9reuphqxtk8gzrrnzd7zqef4qu0ffkvc12]"
$wsTasks.Range("B5").Value = "Take_a_25-minute_walk_without_stopping"
$wsTasks.Range("I5").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("J5").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("L5").Value = " [SECRET, EQUAL, jd5fpzwmdo0mmv]"
$wsTasks.Range("M5").Value = 6
$wsTasks.Range("A6").Value = 3
$wsTasks.Range("B6").Value = "Take_200_steps"
$wsTasks.Range("L6").Value = "[STEPS, STRICTLY_GREATER, 200], [SECRET, EQUAL, tt5wa0bph8wmc0ncswpxv9yxo1zv8nz873jvierr2dg2ta1j]"
$wsTasks.Range("A7").Value = 4
$wsTasks.Range("B7").Value = "tutorial_video(cognitive_activity)"
$wsTasks.Range("F7").Value = "http://localhost:5173/api/media/media-for-ai-b7b4437a/a6cf16fb-1b3c-4862-9086-307cb11c2a41.h5p"
$wsTasks.Range("I7").Value = "H5P_GENERAL"
$wsTasks.Range("J7").Value = "H5P_GENERAL"
$wsTasks.Range("L7").Value = " [SECRET, EQUAL, 98mxf58gtj]"
$wsTasks.Range("A8").Value = 5
$wsTasks.Range("B8").Value = "tutorial_video(social_activity)"
$wsTasks.Range("F8").Value = "http://localhost:5173/api/media/media-for-ai-b7b4437a/b7fb3d01-9712-476f-9d53-4876283c73ce.h5p"
$wsTasks.Range("L8").Value = " [SECRET, EQUAL, nlg7gdk0xywcqlsgk]"
$wsTasks.Range("A9").Value = 6
$wsTasks.Range("B9").Value = "Take_100_steps"
$wsTasks.Range("I9").Value = "WALK"
$wsTasks.Range("J9").Value = "WALK"
$wsTasks.Range("L9").Value = "[STEPS, STRICTLY_GREATER, 100], [SECRET, EQUAL, gtxdio1xesbl]"
$wsTasks.Range("M9").Value = 1
$wsTasks.Range("A10").Value = 6
$wsTasks.Range("B10").Value = "Practice_learning_a_new_skill"
$wsTasks.Range("I10").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("J10").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("L10").Value = " [SECRET, EQUAL, jasbuurj8xir62vluedsnu97oybs]"
$wsTasks.Range("M10").Value = 2
$wsTasks.Range("A11").Value = 6
$wsTasks.Range("B11").Value = "Engage_with_others"
$wsTasks.Range("I11").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("J11").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("L11").Value = " [SECRET, EQUAL, yu1o02lpklafy2yygdso457wcm408tc]"
$wsTasks.Range("B12").Value = "Enjoy_an_activity_with_a_family_member"
$wsTasks.Range("L12").Value = " [SECRET, EQUAL, pys1n3czmhjxfcg9zmp76]"
$wsTasks.Range("M12").Value = 1
$wsTasks.Range("A13").Value = 6
$wsTasks.Range("B13").Value = "Enjoy_an_activity_with_a_family_member"
$wsTasks.Range("F13").ClearContents()
$wsTasks.Range("I13").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("J13").Value = "GENERAL_ACTIVITY"
$wsTasks.Range("L13").Value = " [SECRET, EQUAL, ssiv1kqds4f4s7gbhous]"

# Rows 14-20 (challenge 8 minigame/extra tasks) were removed entirely; delete so the sheet dimension shrinks to A1:N13.
$wsTasks.Range("A14:A20").EntireRow.Delete()
